$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added harvester and experiment design"
# Fill in the new "harvester" values (column B) for all 24 data rows (2-25)
$ws.Range("B2:B25").Value = "S.GISH"

# Fill in the new "experimentDesign" values (column D) for all 24 data rows (2-25)
$ws.Range("D2:D25").Value = "90minuteInduction"

# Fill in the new "strain" values (column F) for all 24 data rows (2-25)
$ws.Range("F2:F25").Value = "KN99allpha"

# Match the resulting selection left by the edit
$ws.Range("F16:F25").Select()
